$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Himanshu Bedi's "Year" cell becomes a free-text note instead of a single year.
$ws.Range("D2").Value = "2020,2022"

# Fix the typo in Himanshu's email address (missing ".ac" before ".in").
$ws.Range("E2").Value = "himanshub.cs.20@nitj.ac.in"

# Turn the three email addresses into live mailto: hyperlinks.
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:himanshub.cs.20@nitj.ac.in")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:ankitj.cs.20@nitj.ac.in")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:harshd.cs.20@nitj.ac.in")

# Resize the columns to fit the new content (values chosen so the saved,
# quantized widths line up with the authored widths as closely as possible).
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 18.0
$ws.Columns.Item(3).ColumnWidth = 16.833333333333332
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 20.5

# Move the active selection to E5, matching where the author left off.
$ws.Range("E5").Select() | Out-Null
